$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newRow = 52

$ws.Cells.Item($newRow, 1).Value = 52
$ws.Cells.Item($newRow, 2).Value = "Login"
$ws.Cells.Item($newRow, 3).Value = "User - jiayu logged in."
$ws.Cells.Item($newRow, 4).Value = "09/05/2022 03:20:46 AM"
